$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Normalized cell-count values (rows 1-23, columns B-G) per the commit
# "Added normalized counts per Grant's request".
$data = @(
    @(1383.1787637855859, 1500.0101265822786, 1638.5331143951835, 1412.6789389435171, 1321.7224390243905, 1225.3431192660551),
    @(1356.1200307771223, 1566.5690174804099, 1831.8800218938152, 1433.2163274639836, 1356.5826829268296, 1279.4568807339449),
    @(1403.8707360861758, 1592.5432188065101, 1845.8075533661743, 1495.5619711868283, 1388.4760975609759, 1311.3321100917431),
    @(1439.6837650679661, 1657.4787221217603, 1924.457142857143, 1572.5771781385777, 1453.7463414634149, 1379.5302752293578),
    @(1520.8599640933573, 1737.0247136829416, 1980.9865353037769, 1648.8589069288819, 1517.5331707317077, 1415.1119266055046),
    @(1621.1364452423697, 1798.7134418324292, 2094.8645867542423, 1730.2749828493027, 1601.346097560976, 1518.1504587155962),
    @(1675.2539112592972, 1875.0126582278483, 2146.4783798576905, 1794.8210610564831, 1694.8012195121955, 1587.8311926605504),
    @(1749.2675044883304, 1937.5130801687765, 2230.0435686918445, 1880.6380059455753, 1742.2704878048785, 1650.0990825688073),
    @(1835.2187740446268, 2022.7409282700423, 2302.9582922824302, 1926.8471301166248, 1817.1829268292688, 1733.8642201834862),
    @(1920.3741985124391, 2139.6248342374925, 2434.0409414340452, 2057.4062428538764, 1885.4200000000005, 1785.7541284403669),
    @(2089.8892023595795, 2330.3728752260399, 2632.3034482758621, 2193.0997027212443, 2070.1051219512201, 1952.5431192660549),
    @(2221.203641959477, 2422.9059674502714, 2783.8677613574168, 2375.7357649211071, 2219.9300000000007, 2054.0990825688073),
    @(2368.4349833290585, 2576.3160940325502, 2956.7330049261086, 2503.3609650125773, 2331.9278048780493, 2191.2366972477062),
    @(2526.0123108489356, 2825.5060880048222, 3172.2001094690754, 2672.7944203064258, 2497.3285365853662, 2316.5137614678897),
    @(2681.2021031033596, 2935.0847498493072, 3313.1139573070609, 2785.7500571689916, 2623.4187804878056, 2457.3577981651374),
    @(2827.6375993844576, 3064.1440626883668, 3407.3296113847841, 2927.3113423279215, 2799.9451219512202, 2533.7100917431189),
    @(2907.2221082328801, 3173.7227245328513, 3557.2553913519432, 3062.2713240338444, 2849.6395121951227, 2638.9724770642201),
    @(2973.2772505770708, 3260.5739602169983, 3613.7847837985773, 3087.2095815229823, 2943.0946341463423, 2702.7229357798165),
    @(3014.6611951782506, 3293.0417118746236, 3682.6031746031749, 3134.8856620169222, 3000.9478048780493, 2765.7321100917429),
    @(3040.9240830982303, 3341.7433393610613, 3709.6389709906953, 3133.4187056940318, 3009.1065853658542, 2780.5577981651372),
    @(3080.7163375224418, 3322.2626883664861, 3701.4463054187195, 3142.95392179282, 3017.2653658536592, 2826.5174311926603),
    @(3100.6124647345473, 3340.9316455696207, 3703.9041050903124, 3143.6873999542649, 3058.8009756097567, 2814.6568807339449),
    @(3102.2041549115156, 3357.9772151898737, 3729.3013683634376, 3186.2291333180883, 3046.9336585365859, 2823.5522935779813),
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 1
    for ($j = 0; $j -lt $data[$i].Length; $j++) {
        $col = $j + 2
        $ws.Cells.Item($row, $col).Value = $data[$i][$j]
    }
}

# Restore selection as left by the author after the edit
[void]$ws.Range("B1:G23").Select()
